$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New row 12 - V6 / Sliding window
$ws.Range("A12").Value = "V6"
$ws.Range("B12").Value = "Sliding window"

# Add value to row 7 (existing row) - "1941, 7, feature_size"
$ws.Range("C7").Value = "1941, 7, feature_size"

# New row 13 - V7 / Without sliding window / 1, 1941, feature_size
$ws.Range("A13").Value = "V7"
$ws.Range("B13").Value = "Without sliding window"
$ws.Range("C13").Value = "1, 1941, feature_size"

# Update the selection to match the final cursor position seen in the diff
$ws.Range("C14").Select()
